$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (was "leaf.area", now "HL_cover")
$ws.Name = "HL_cover"

# Row 2: (Intercept) - update Estimate / Std.Error / t value (Pr(>|t|) stays 0)
$ws.Range("B2").Value = 90.767
$ws.Range("C2").Value = 9.303
$ws.Range("D2").Value = 9.756
$ws.Range("E2").Value = 0

# Row 3: renamed from soil_cover to management2entbuscht, with new stats
$ws.Range("A3").Value = "management2entbuscht"
$ws.Range("B3").Value = -12.613
$ws.Range("C3").Value = 5.588
$ws.Range("D3").Value = -2.257
$ws.Range("E3").Value = 0.027

# New row 4: management2Buche
$ws.Range("A4").Value = "management2Buche"
$ws.Range("B4").Value = -9.471
$ws.Range("C4").Value = 5.771
$ws.Range("D4").Value = -1.641
$ws.Range("E4").Value = 0.1051

# New row 5: management2Fichte
$ws.Range("A5").Value = "management2Fichte"
$ws.Range("B5").Value = -35.861
$ws.Range("C5").Value = 6.39
$ws.Range("D5").Value = -5.612
$ws.Range("E5").Value = 0

# New row 6: SL_cover
$ws.Range("A6").Value = "SL_cover"
$ws.Range("B6").Value = -0.202
$ws.Range("C6").Value = 0.09
$ws.Range("D6").Value = -2.245
$ws.Range("E6").Value = 0.0278

# New row 7: soil_cover
$ws.Range("A7").Value = "soil_cover"
$ws.Range("B7").Value = -0.564
$ws.Range("C7").Value = 0.12
$ws.Range("D7").Value = -4.689
$ws.Range("E7").Value = 0

# New row 8: moss_cover
$ws.Range("A8").Value = "moss_cover"
$ws.Range("B8").Value = -0.247
$ws.Range("C8").Value = 0.1
$ws.Range("D8").Value = -2.473
$ws.Range("E8").Value = 0.0157

# Apply the same style (border etc.) as the existing data rows to the new rows
$srcStyle = $ws.Range("A3:E3")
$dstStyle = $ws.Range("A4:E8")
$srcStyle.Copy()
$dstStyle.PasteSpecial(-4122) # xlPasteFormats
